$d = $word.ActiveDocument

# Locate the anchor paragraph - the sentence right before the blank line
# that should receive the new date_range merge field.
$anchorIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*booked out to your section*") {
        $anchorIdx = $i
        break
    }
}

if ($anchorIdx -eq -1) {
    throw "Could not locate anchor paragraph"
}

# The target is the second blank paragraph following the anchor
# (there is a run of blank paragraphs between the notification
# sentence and the "Please copy all records..." instruction).
$blankSeen = 0
$targetIdx = -1
for ($i = $anchorIdx + 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "") {
        $blankSeen = $blankSeen + 1
        if ($blankSeen -eq 2) {
            $targetIdx = $i
            break
        }
    } else {
        break
    }
}

if ($targetIdx -eq -1) {
    throw "Could not locate target paragraph"
}

$targetRange = $d.Paragraphs.Item($targetIdx).Range

$newParaXml = @'
<w:p w14:paraId="6C7AD2C3" w14:textId="77777777" w:rsidR="000E5C13" w:rsidRDefault="000E5C13" w:rsidP="00DB25D8"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="24"/><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="24"/><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="24"/><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =date_range \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="24"/><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:noProof/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="24"/><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:t>«=date_range»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="24"/><w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr"><w14:noFill/><w14:prstDash w14:val="solid"/><w14:round/></w14:textOutline></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@

[void]$targetRange.InsertXML($newParaXml)

Write-Host "Updated paragraph $targetIdx with date_range merge field"
